$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Coin(2), C=Link(3), D=Price(4), E=Volume(1h)(5)

# Helper: write a value as TEXT (preventing Excel's automatic
# string->number coercion for numeric-looking strings), while
# preserving the cell's original style.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "35.679.97"
Set-TextValue 2 5 "  +3.33%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "1.859.26"
Set-TextValue 3 5 "  +2.45%  "

# Row 5 - BNB
Set-TextValue 5 4 "230.79"
Set-TextValue 5 5 "  +1.90%  "

# Row 6 - XRP
Set-TextValue 6 4 "0.613"
Set-TextValue 6 5 "  +2.98%  "

# Row 7 - USDC
Set-TextValue 7 5 "  +0.22%  "

# Row 8 - Solana
Set-TextValue 8 4 "42.78"
Set-TextValue 8 5 "  +11.90%  "

# Row 9 - Cardano
Set-TextValue 9 5 "  +6.78%  "

# Row 10 - Dogecoin
Set-TextValue 10 5 "  +2.50%  "

# Row 11 - TRON
Set-TextValue 11 4 "0.101"
Set-TextValue 11 5 "  +2.89%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue 12 4 "2.128.76"
Set-TextValue 12 5 "  +2.51%  "

# Row 13 - Chainlink
Set-TextValue 13 4 "11.53"
Set-TextValue 13 5 "  +2.03%  "

# Row 14 - WrappedEther
Set-TextValue 14 4 "1.861.73"
Set-TextValue 14 5 "  +2.22%  "

# Row 15 - Polygon
Set-TextValue 15 4 "0.681"
Set-TextValue 15 5 "  +7.04%  "

# Row 16 - Polkadot
Set-TextValue 16 4 "4.78"
Set-TextValue 16 5 "  +7.51%  "

# Row 17 - WrappedBTC
Set-TextValue 17 4 "35.679.89"
Set-TextValue 17 5 "  +3.40%  "

# Row 18 - Litecoin
Set-TextValue 18 4 "70.24"
Set-TextValue 18 5 "  +2.28%  "

# Row 19 - ShibaInu
Set-TextValue 19 4 "0.0₃0799"
Set-TextValue 19 5 "  +2.78%  "

# Row 20 - BitcoinCash
Set-TextValue 20 4 "246.03"
Set-TextValue 20 5 "  +0.82%  "

# Row 21 - Avalanche
Set-TextValue 21 4 "12.21"
Set-TextValue 21 5 "  +8.04%  "

# Row 22 - Uniswap
Set-TextValue 22 5 "  +14.63%  "

# Row 23 - Dai
Set-TextValue 23 5 "  +0.18%  "

# Row 24 - Toncoin
Set-TextValue 24 4 "2.23"
Set-TextValue 24 5 "  +0.16%  "

# Row 25 - Monero
Set-TextValue 25 4 "170.54"
Set-TextValue 25 5 "  -0.02%  "

# Row 26 - Cosmos
Set-TextValue 26 4 "7.97"
Set-TextValue 26 5 "  +0.10%  "

# Row 27 - EthereumClassic
Set-TextValue 27 4 "17.91"
Set-TextValue 27 5 "  +0.77%  "

# Row 28 - Stellar
Set-TextValue 28 4 "0.124"
Set-TextValue 28 5 "  +1.84%  "

# Row 29 - PancakeSwap
Set-TextValue 29 5 "  +15.94%  "

# Row 30 - BinanceUSD
Set-TextValue 30 4 "1.00"
Set-TextValue 30 5 "  +0.14%  "

# Row 31 - EURNeutrino
Set-TextValue 31 4 "3.327.21"
Set-TextValue 31 5 "  +36.94%  "

# Row 32 - was Hedera, becomes Filecoin
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 32 4 "3.93"
Set-TextValue 32 5 "  +3.17%  "

# Row 33 - was Filecoin, becomes Hedera
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 33 4 "0.0543"
Set-TextValue 33 5 "  +4.14%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue 34 5 "  +4.62%  "

# Row 35 - LidoDAOToken
Set-TextValue 35 5 "  +2.86%  "

# Row 36 - ImmutableX
Set-TextValue 36 4 "0.682"
Set-TextValue 36 5 "  +4.97%  "

# Row 37 - RenderToken
Set-TextValue 37 4 "2.54"
Set-TextValue 37 5 "  +7.52%  "

# Row 38 - was TrustWalletToken, becomes Aave
$ws.Cells.Item(38, 2).Value = "Aave"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 38 4 "88.98"
Set-TextValue 38 5 "  +8.65%  "

# Row 39 - was Aave, becomes TrustWalletToken
$ws.Cells.Item(39, 2).Value = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue 39 4 "1.09"
Set-TextValue 39 5 "  +2.57%  "

# Row 40 - Maker
Set-TextValue 40 4 "1.342.39"
Set-TextValue 40 5 "  -1.98%  "

# Row 41 - VeChain
Set-TextValue 41 5 "  +4.35%  "

# Row 42 - ARBITRUM
Set-TextValue 42 4 "1.02"
Set-TextValue 42 5 "  +7.00%  "

# Row 43 - WEMIXToken
Set-TextValue 43 5 "  +5.61%  "

# Row 44 - InjectiveProtocol
Set-TextValue 44 4 "15.14"
Set-TextValue 44 5 "  +8.82%  "

# Row 45 - HuobiToken
Set-TextValue 45 4 "2.48"
Set-TextValue 45 5 "  +1.75%  "

# Row 46 - MXToken
Set-TextValue 46 4 "2.82"
Set-TextValue 46 5 "  +1.37%  "

# Row 47 - Kaspa
Set-TextValue 47 5 "  +2.66%  "

# Row 48 - FraxShare
Set-TextValue 48 4 "6.11"
Set-TextValue 48 5 "  +5.04%  "

# Row 49 - RocketPoolETH
Set-TextValue 49 4 "2.026.67"
Set-TextValue 49 5 "  +2.48%  "

# Row 50 - Quant
Set-TextValue 50 4 "104.54"
Set-TextValue 50 5 "  +2.09%  "

# Row 51 - PaxDollar
Set-TextValue 51 5 "  +0.20%  "
